$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# A1 holds a date serial; bump it by one day (45310 -> 45311, i.e. 2024-01-19 -> 2024-01-20)
$ws.Range("A1").Value = Get-Date -Year 2024 -Month 1 -Day 20 -Hour 0 -Minute 0 -Second 0

# D19 price update (1051 -> 440)
$ws.Range("D19").Value = 440
